# Apply the "rozbor prace" workbook update:
# - rename the last category from "dokumentace" to "prezentace"
# - update several hour counters in column H
# - add the missing H22 value and extend the H23 SUM formula to include it
# - move the active cell selection to I23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C22: dokumentace -> prezentace
$ws.Range("C22").Value = "prezentace"

# Column H hour adjustments
$ws.Range("H8").Value = 4
$ws.Range("H9").Value = 5
$ws.Range("H11").Value = 3
$ws.Range("H13").Value = 4
$ws.Range("H14").Value = 8
$ws.Range("H15").Value = 3
$ws.Range("H22").Value = 2

# Extend the total formula to cover the new H22 value
$ws.Range("H23").Formula = "=SUM(H2:H22)"

# Move the selected cell, matching the saved view state
$ws.Range("I23").Select()
